# Apply corrected Diebold-Mariano test statistics (C) and P-values (D)
# for rows 2-11 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.09895831063788695
$ws.Range("D2").Value = 0.9220675324683478

$ws.Range("C3").Value = 0.7804417311379241
$ws.Range("D3").Value = 0.4434476219107353

$ws.Range("C4").Value = 0.8116873613062648
$ws.Range("D4").Value = 0.4256665352610445

$ws.Range("C5").Value = -0.3037694819639957
$ws.Range("D5").Value = 0.7641578843725607

$ws.Range("C6").Value = 0.5023987842197064
$ws.Range("D6").Value = 0.6203788460323709

$ws.Range("C7").Value = 0.6177294781349821
$ws.Range("D7").Value = 0.5430957605406712

$ws.Range("C8").Value = -0.3496547526933666
$ws.Range("D8").Value = 0.7299227848561594

$ws.Range("C9").Value = 0.1311459260075785
$ws.Range("D9").Value = 0.8968517774632179

$ws.Range("C10").Value = -0.8385523144391922
$ws.Range("D10").Value = 0.4107394865044696

$ws.Range("C11").Value = -0.8961006664061658
$ws.Range("D11").Value = 0.3798998951930821
